$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Platform Coverage")
$ws2 = $wb.Worksheets.Item("MarketShare")

# --- Data edit: make the Trichuris scenario-1 coverage yearly instead of
# biennial by filling in the odd "skipped" years (I, K, M, O, Q, S, U, W,
# Y, AA, AC) in row 2 of "Platform Coverage" with the same 0.6 coverage
# value already present in the even years.
$oddYearCols = @("I", "K", "M", "O", "Q", "S", "U", "W", "Y", "AA", "AC")
foreach ($col in $oddYearCols) {
    $ws1.Range($col + "2").Value = 0.6
}

# --- View-state: restore each sheet's on-screen selection/zoom, leaving
# "MarketShare" as the last-activated (and therefore active) sheet.
$ws1.Select()
$ws1.Range("AE2").Select()
$excel.ActiveWindow.Zoom = 166

$ws2.Select()
$ws2.Range("Z3").Select()
$excel.ActiveWindow.Zoom = 181
